$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "radio" value with "choice" in the type column (C2:C45)
$range = $ws.Range("C2:C45")
$range.Value = "choice"

# Select the edited range and scroll the frozen-pane view down,
# matching the end-state view captured after the edit.
$ws.Range("A34").Select()
$range.Select()
